$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.468.20'
$ws.Range('D3').Value = '1.560.57'
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('E4').Value = '  +0.72%  '
$ws.Range('D5').Value = '''211.58'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('E6').Value = '  -0.65%  '
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('D8').Value = '''45.96'
$ws.Range('E8').Value = '  +3.52%  '
$ws.Range('D9').Value = '''24.05'
$ws.Range('E9').Value = '  +0.45%  '
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('E12').Value = '  -0.55%  '
$ws.Range('D13').Value = '1.784.42'
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').Value = '1.583.67'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('E15').Value = '  -2.40%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '28.476.04'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = '''3.67'
$ws.Range('E17').Value = '  -2.91%  '
$ws.Range('D18').Value = '''61.84'
$ws.Range('E18').Value = '  -3.36%  '
$ws.Range('D19').Value = '''226.58'
$ws.Range('E19').Value = '  -3.67%  '
$ws.Range('D20').Value = '0.0₃0693'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('E21').Value = '  -2.88%  '
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '''3.87'
$ws.Range('E23').Value = '  -6.66%  '
$ws.Range('E24').Value = '  -3.37%  '
$ws.Range('D25').Value = '''2.09'
$ws.Range('E25').Value = '  +7.47%  '
$ws.Range('D26').Value = '''149.79'
$ws.Range('E26').Value = '  -1.10%  '
$ws.Range('D27').Value = '''14.94'
$ws.Range('E27').Value = '  -2.73%  '
$ws.Range('E28').Value = '  -3.12%  '
$ws.Range('E29').Value = '  -3.00%  '
$ws.Range('E30').Value = '  +0.68%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.0464'
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '''1.11'
$ws.Range('E32').Value = '  -3.71%  '
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('D35').Value = '1.394.65'
$ws.Range('E35').Value = '  -1.79%  '
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('E37').Value = '  -4.36%  '
$ws.Range('E38').Value = '  +1.84%  '
$ws.Range('E39').Value = '  -0.42%  '
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = '''0.534'
$ws.Range('E41').Value = '  -2.10%  '
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').Value = '''0.786'
$ws.Range('E43').Value = '  -3.57%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''1.85'
$ws.Range('E44').Value = '  +0.74%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '''5.54'
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('D46').Value = '''0.979'
$ws.Range('E46').Value = '  +0.87%  '
$ws.Range('D47').Value = '''62.69'
$ws.Range('E47').Value = '  -2.89%  '
$ws.Range('D48').Value = '1.697.55'
$ws.Range('E48').Value = '  -1.67%  '
$ws.Range('D49').Value = '''85.89'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('E50').Value = '  +2.44%  '
$ws.Range('E51').Value = '  -1.39%  '
